$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text storage for the Price/Volume columns so values like
# "30.456.64" (multi-dot) and plain decimals stay text, matching the
# original inlineStr cells instead of being reinterpreted as numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "30.456.64"
$ws.Range("E2").Value = "  +1.58%  "
$ws.Range("D3").Value = "1.998.91"
$ws.Range("E3").Value = "  +4.12%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "324.26"
$ws.Range("E5").Value = "  +1.25%  "
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  +1.44%  "
$ws.Range("D8").Value = "0.4144"
$ws.Range("E8").Value = "  +2.99%  "
$ws.Range("D9").Value = "0.08732"
$ws.Range("E9").Value = "  +5.93%  "
$ws.Range("D10").Value = "1.134"
$ws.Range("E10").Value = "  +2.08%  "
$ws.Range("D11").Value = "42.98"
$ws.Range("E11").Value = "  +2.25%  "
$ws.Range("E12").Value = "  +4.74%  "
$ws.Range("D13").Value = "1.994.36"
$ws.Range("E13").Value = "  +4.26%  "
$ws.Range("D14").Value = "6.584"
$ws.Range("E14").Value = "  +2.86%  "
$ws.Range("D15").Value = "7.438"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "94.19"
$ws.Range("E17").Value = "  +2.21%  "
$ws.Range("D18").Value = "0.00001118"
$ws.Range("E18").Value = "  +2.12%  "
$ws.Range("D19").Value = "0.06512"
$ws.Range("E19").Value = "  +0.36%  "
$ws.Range("D20").Value = "18.97"
$ws.Range("E20").Value = "  +4.86%  "
$ws.Range("E21").Value = "  -0.07%  "
$ws.Range("D22").Value = "6.174"
$ws.Range("E22").Value = "  +3.74%  "
$ws.Range("D23").Value = "30.515.97"
$ws.Range("E23").Value = "  +1.61%  "
$ws.Range("E24").Value = "  +4.91%  "
$ws.Range("D25").Value = "2.228"
$ws.Range("E25").Value = "  +1.79%  "
$ws.Range("D26").Value = "2.227.19"
$ws.Range("E26").Value = "  +4.31%  "
$ws.Range("D27").Value = "22.37"
$ws.Range("E27").Value = "  +0.31%  "
$ws.Range("D28").Value = "163.34"
$ws.Range("E28").Value = "  +1.17%  "
$ws.Range("D29").Value = "2.405"
$ws.Range("E29").Value = "  +3.49%  "
$ws.Range("D30").Value = "131.41"
$ws.Range("E30").Value = "  +1.99%  "
$ws.Range("D31").Value = "1.138"
$ws.Range("E31").Value = "  +0.92%  "
$ws.Range("E32").Value = "  +0.90%  "
$ws.Range("D33").Value = "6.084"
$ws.Range("E33").Value = "  +1.65%  "
$ws.Range("D34").Value = "3.854"
$ws.Range("E34").Value = "  +0.88%  "
$ws.Range("D35").Value = "1.334"
$ws.Range("E35").Value = "  +11.43%  "
$ws.Range("E36").Value = "  +2.95%  "
$ws.Range("D37").Value = "5.433"
$ws.Range("E37").Value = "  +0.86%  "
$ws.Range("D38").Value = "0.06596"
$ws.Range("D39").Value = "12.40"
$ws.Range("E39").Value = "  +9.00%  "
$ws.Range("E40").Value = "  +2.16%  "
$ws.Range("D41").Value = "9.029"
$ws.Range("E41").Value = "  +0.88%  "
$ws.Range("D42").Value = "0.6621"
$ws.Range("E42").Value = "  +3.14%  "
$ws.Range("D43").Value = "1.238"
$ws.Range("E43").Value = "  +1.79%  "
$ws.Range("D44").Value = "13.57"
$ws.Range("E44").Value = "  +1.35%  "
$ws.Range("D45").Value = "0.6162"
$ws.Range("E45").Value = "  +2.72%  "
$ws.Range("E46").Value = "  +1.33%  "
$ws.Range("D47").Value = "3.663"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("D48").Value = "1.270"
$ws.Range("E48").Value = "  +4.54%  "
$ws.Range("D49").Value = "124.32"
$ws.Range("E49").Value = "  +1.19%  "
$ws.Range("D50").Value = "80.00"
$ws.Range("E50").Value = "  +1.56%  "
$ws.Range("D51").Value = "0.06898"
$ws.Range("E51").Value = "  +1.59%  "

# Restore the default (unstyled) cell style so only the values change,
# matching the original workbook which had no explicit style on these cells.
$ws.Range("D2:E51").Style = "Normal"
